$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95. This shifts the existing rows 95-319
# down to 96-320 (each row keeps its own original data), matching the
# target diff which shows every record from row 95 onward moving down
# by one row, with a brand-new record introduced at row 95 and the
# previously-last record (old row 319) now living in the new row 320.
$ws.Rows(95).Insert()

# Populate the newly inserted row 95 with the new record's data.
$ws.Cells.Item(95, 1).Value = 3
$ws.Cells.Item(95, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(95, 3).Value = "Coquimbo"
$ws.Cells.Item(95, 4).Value = 44708
$ws.Cells.Item(95, 5).Value = 5
$ws.Cells.Item(95, 6).Value = 100112039
$ws.Cells.Item(95, 7).Value = "Ciboulette"
$ws.Cells.Item(95, 8).Value = "Sin especificar"
$ws.Cells.Item(95, 9).Value = "Primera"
$ws.Cells.Item(95, 10).Value = 160
$ws.Cells.Item(95, 11).Value = 1500
$ws.Cells.Item(95, 12).Value = 1500
$ws.Cells.Item(95, 13).Value = 1500
$ws.Cells.Item(95, 14).Value = "$/docena de atados"
$ws.Cells.Item(95, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(95, 16).Value = 500
$ws.Cells.Item(95, 17).Value = 3
$ws.Cells.Item(95, 18).Value = "Hortaliza"
